$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "statut" column (A) used four emoji values as status codes.
# Replace them with the new textual/emoji labels:
#   old -> new
#   📕 (red book)    -> -3
#   📘 (blue book)   -> ⚠️
#   📙 (orange book) -> +3
#   📗 (green book)  -> ✅

$map = @{
    "📕" = "-3"
    "📘" = "⚠️"
    "📙" = "+3"
    "📗" = "✅"
}

$used = $ws.UsedRange
$lastRow = $used.Rows.Count + $used.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $old = $cell.Value2
    if ($map.ContainsKey($old)) {
        $new = $map[$old]
        $cell.NumberFormat = "@"
        $cell.Value = $new
    }
}
